$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the outlier data row (Total Width > 1000) -- row 146, the
# launchEditor.js entry (Total Width = 7287) under
# packages/react-dev-utils. Deleting the entire row shifts everything
# below it up by one and Excel auto-adjusts the MAX/MIN/AVERAGE/STDEV.P
# formula ranges (C2:C252 -> C2:C251, etc.) and the sheet dimension.
$ws.Rows(146).Delete()

# Restore the selection to match the saved state after the edit.
$ws.Range("O11").Select()
